$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 94.73967398981497
$ws.Range("C2").Value = 130.5725205260365
$ws.Range("D2").Value = 148.3532478954353
$ws.Range("E2").Value = 152.5655404845378

$ws.Range("B3").Value = 115.8526779330273
$ws.Range("C3").Value = 156.8049429565898
$ws.Range("D3").Value = 175.3102001278824
$ws.Range("E3").Value = 181.2896197755404

$ws.Range("B4").Value = 116.8445518566937
$ws.Range("C4").Value = 157.7018885027459
$ws.Range("D4").Value = 175.2974291660801
$ws.Range("E4").Value = 180.9409395364281

$ws.Range("B5").Value = 92.98295278191156
$ws.Range("C5").Value = 116.6419887944614
$ws.Range("D5").Value = 123.4614096019962
$ws.Range("E5").Value = 127.3609302410908

$ws.Range("B6").Value = 80.73943360570159
$ws.Range("C6").Value = 102.1626600216092
$ws.Range("D6").Value = 109.3217471349291
$ws.Range("E6").Value = 111.4236309222504

$ws.Range("B7").Value = 9.333911543418484
$ws.Range("C7").Value = 11.23422062765218
$ws.Range("D7").Value = 11.78526578965354
$ws.Range("E7").Value = 11.8391811533395

$ws.Range("B8").Value = 205.1106828654895
$ws.Range("C8").Value = 378.8514805626896
$ws.Range("D8").Value = 489.5722118061616
$ws.Range("E8").Value = 583.5252692555785

$ws.Range("B9").Value = 108.8989953560731
$ws.Range("C9").Value = 140.3493981218891
$ws.Range("D9").Value = 153.5434512669169
$ws.Range("E9").Value = 158.3777453323874

$ws.Range("B10").Value = 51.8584832674242
$ws.Range("C10").Value = 63.38079649864315
$ws.Range("D10").Value = 67.99010038631985
$ws.Range("E10").Value = 67.62793519753488

$ws.Range("B11").Value = 9.584581920080701
$ws.Range("C11").Value = 11.08494337666994
$ws.Range("D11").Value = 11.74222911804124
$ws.Range("E11").Value = 12.60110022955488

$ws.Range("B12").Value = 23.83717678683425
$ws.Range("C12").Value = 29.85737786738835
$ws.Range("D12").Value = 31.42426284019058
$ws.Range("E12").Value = 30.75470296023732

$ws.Range("B13").Value = 29.13133318340149
$ws.Range("C13").Value = 35.31214364868075
$ws.Range("D13").Value = 38.24731927604619
$ws.Range("E13").Value = 38.4931219212982
